# This script updates the test fixture data in before.xlsx so that the
# "22" batch of generated users becomes the "23" batch, and the "27" batch
# becomes the "28" batch, matching the author's commit.
#
# Sheet1 ("Sheet1") holds the canonical generator: column I has the numeric
# suffix (I2 = 22, I23 = 27) and columns A/B/C hold formulas
# (=CONCATENATE(...)) that derive the user name / email from it, so Excel
# recalculates those automatically once I2 / I23 change.
#
# The "login" and "order" sheets hold literal (non-formula) copies of the
# generated values in columns G/H/I and R/S/T respectively, so those need
# to be updated explicitly to keep them in sync.

$wb = $excel.ActiveWorkbook

$loginSheet = $wb.Worksheets.Item("login")
$orderSheet = $wb.Worksheets.Item("order")
$dataSheet  = $wb.Worksheets.Item("Sheet1")

# --- Sheet1: bump the two numeric generator seeds ---------------------
$dataSheet.Range("I2").Value = 23
$dataSheet.Range("I23").Value = 28

# --- login sheet: literal new_user_name / new_email values ------------
$loginNames = @(
    "EthanBaker", "DelanieCarman", "BretAgnew", "EdgardoTaylor", "TyrekReis",
    "LeannaChow", "TuckerCarlson", "AnnmarieConnor", "MoniqueWitte",
    "MikelWhitlock", "VincentAmaya", "KeiraQuiroz", "EllisCreech",
    "DionteCreel", "NicholeFoust", "ManuelConnell", "LourdesElam",
    "LincolnFrederick", "AlisaCash", "LucilleGriffiths"
)

for ($i = 0; $i -lt $loginNames.Count; $i++) {
    $row = 2 + $i
    $newName = "{0}23" -f $loginNames[$i]
    $newEmail = "{0}23@gmail.com" -f $loginNames[$i]

    $loginSheet.Range("G$row").Value = $newName
    $loginSheet.Range("H$row").Value = $newName
    $loginSheet.Range("I$row").Value = $newEmail
}

# --- order sheet: literal new_user_name / new_email values -------------
$orderNames = @(
    "DonnellJernigan", "MalikOtoole", "AlanCaudill", "AdanApplegate",
    "AiyanaWhitworth", "MercedezBrien", "DuaneHager", "LorenBell",
    "GeraldHiller", "DeionBranch", "DakotaHalstead", "ElliottFurman",
    "MiltonCamp", "DawnChester", "ZacheryPetrie", "EstebanAngel",
    "JimmyBlankenship", "AllysaGrice", "AugustineYoo", "BrandiSouthard"
)

for ($i = 0; $i -lt $orderNames.Count; $i++) {
    $row = 2 + $i
    $newName = "{0}28" -f $orderNames[$i]
    $newEmail = "{0}28@gmail.com" -f $orderNames[$i]

    $orderSheet.Range("R$row").Value = $newName
    $orderSheet.Range("S$row").Value = $newName
    $orderSheet.Range("T$row").Value = $newEmail
}
